$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 22.320355
$ws.Range("H2").Value = 66.961065
$ws.Range("I2").Value = 0.4795709034536316
$ws.Range("J2").Value = 0.480894303533459
$ws.Range("M2").Value = 0.1938303333333333
$ws.Range("N2").Value = 0.581491
$ws.Range("Q2").Value = 4.326361849768333
$ws.Range("R2").Value = 38.937256647915
$ws.Range("S2").Value = 0.4795709034536316
$ws.Range("T2").Value = 0.480894303533459

# Row 3
$ws.Range("I3").Value = 0.003471397841042084
$ws.Range("J3").Value = 0.003480977338352835
$ws.Range("M3").Value = 0.1938303333333333
$ws.Range("N3").Value = 0.581491
$ws.Range("Q3").Value = 0.03131658546566667
$ws.Range("R3").Value = 0.281849269191
$ws.Range("S3").Value = 0.003471397841042084
$ws.Range("T3").Value = 0.003480977338352835

# Row 4
$ws.Range("G4").Value = 16.184686
$ws.Range("H4").Value = 48.554058
$ws.Range("I4").Value = 0.3477410859788449
$ws.Range("J4").Value = 0.3487006950327504
$ws.Range("M4").Value = 0.1938303333333333
$ws.Range("N4").Value = 0.581491
$ws.Range("Q4").Value = 3.137083082275333
$ws.Range("R4").Value = 28.233747740478
$ws.Range("S4").Value = 0.3477410859788449
$ws.Range("T4").Value = 0.3487006950327504

# Row 5
$ws.Range("G5").Value = 0.3842475
$ws.Range("H5").Value = 0.7684949999999999
$ws.Range("I5").Value = 0.008255868722733095
$ws.Range("J5").Value = 0.005519100805728606
$ws.Range("M5").Value = 0.1938303333333333
$ws.Range("N5").Value = 0.581491
$ws.Range("Q5").Value = 0.07447882100749999
$ws.Range("R5").Value = 0.446872926045
$ws.Range("S5").Value = 0.008255868722733095
$ws.Range("T5").Value = 0.005519100805728606

# Row 6
$ws.Range("G6").Value = 7.491490666666667
$ws.Range("H6").Value = 22.474472
$ws.Range("I6").Value = 0.1609607440037482
$ws.Range("J6").Value = 0.1614049232897091
$ws.Range("M6").Value = 0.1938303333333333
$ws.Range("N6").Value = 0.581491
$ws.Range("Q6").Value = 1.452078133083556
$ws.Range("R6").Value = 13.068703197752
$ws.Range("S6").Value = 0.1609607440037482
$ws.Range("T6").Value = 0.1614049232897091
